# Remove the 2000/2005/2006/2007/2008/2009 data rows (old rows 2-7),
# shifting the 2010-2013 rows up so they become rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A7").EntireRow.Delete()
